$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.127786734538432
$ws.Cells.Item(2, 4).Value = 0.2787545109805052
$ws.Cells.Item(2, 5).Value = 0.1661596229618656
$ws.Cells.Item(2, 6).Value = 4.346183224801223
$ws.Cells.Item(2, 7).Value = 0.002580812038164213
$ws.Cells.Item(2, 11).Value = 0.4603832145521665
$ws.Cells.Item(2, 12).Value = 0.1622099368341594
$ws.Cells.Item(2, 13).Value = 0.2446368497948157
$ws.Cells.Item(3, 2).Value = 1.115670996603114
$ws.Cells.Item(3, 4).Value = 0.2670270745800849
$ws.Cells.Item(3, 5).Value = 0.16377099887942
$ws.Cells.Item(3, 6).Value = 4.136223002071262
$ws.Cells.Item(3, 7).Value = 0.002586348066413418
$ws.Cells.Item(3, 11).Value = 0.4312135845364367
$ws.Cells.Item(3, 12).Value = 0.1589314389476471
$ws.Cells.Item(3, 13).Value = 0.2414107227606124
$ws.Cells.Item(4, 2).Value = 1.10888833901825
$ws.Cells.Item(4, 4).Value = 0.2597694320987216
$ws.Cells.Item(4, 5).Value = 0.1622515330430194
$ws.Cells.Item(4, 6).Value = 4.007728903964875
$ws.Cells.Item(4, 7).Value = 0.002589921689355165
$ws.Cells.Item(4, 11).Value = 0.4140218462186738
$ws.Cells.Item(4, 12).Value = 0.1570602074369702
$ws.Cells.Item(4, 13).Value = 0.2395949333201912
$ws.Cells.Item(5, 2).Value = 1.106289480948675
$ws.Cells.Item(5, 4).Value = 0.2567971514159808
$ws.Cells.Item(5, 5).Value = 0.1616188763680442
$ws.Cells.Item(5, 6).Value = 3.955467725124038
$ws.Cells.Item(5, 7).Value = 0.002591422003218685
$ws.Cells.Item(5, 11).Value = 0.4071950506422866
$ws.Cells.Item(5, 12).Value = 0.1563332267837154
$ws.Cells.Item(5, 13).Value = 0.2388964702236791
$ws.Cells.Item(6, 2).Value = 1.105867917916726
$ws.Cells.Item(6, 4).Value = 0.2563027046593049
$ws.Cells.Item(6, 5).Value = 0.1615130060444718
$ws.Cells.Item(6, 6).Value = 3.946795764321791
$ws.Cells.Item(6, 7).Value = 0.002591673793269337
$ws.Cells.Item(6, 11).Value = 0.4060722268737038
$ws.Cells.Item(6, 12).Value = 0.1562146574358394
$ws.Cells.Item(6, 13).Value = 0.2387829960332404
$ws.Cells.Item(7, 2).Value = 1.108852621224372
$ws.Cells.Item(7, 4).Value = 0.2597294070689742
$ws.Cells.Item(7, 5).Value = 0.1622430555472292
$ws.Cells.Item(7, 6).Value = 4.007023688219761
$ws.Cells.Item(7, 7).Value = 0.002589941744642711
$ws.Cells.Item(7, 11).Value = 0.4139290552354282
$ws.Cells.Item(7, 12).Value = 0.1570502592531966
$ws.Cells.Item(7, 13).Value = 0.2395853456555237
$ws.Cells.Item(8, 2).Value = 1.123473006905613
$ws.Cells.Item(8, 4).Value = 0.2747224594085793
$ws.Cells.Item(8, 5).Value = 0.1653468868667298
$ws.Cells.Item(8, 6).Value = 4.27369820781874
$ws.Cells.Item(8, 7).Value = 0.002582684743755371
$ws.Cells.Item(8, 11).Value = 0.4501753518582916
$ws.Cells.Item(8, 12).Value = 0.1610500277635154
$ws.Cells.Item(8, 13).Value = 0.2434902065329858
$ws.Cells.Item(9, 2).Value = 1.157353167915971
$ws.Cells.Item(9, 4).Value = 0.3036916516762744
$ws.Cells.Item(9, 5).Value = 0.1710222007816533
$ws.Cells.Item(9, 6).Value = 4.800248952676895
$ws.Cells.Item(9, 7).Value = 0.002569830910348838
$ws.Cells.Item(9, 11).Value = 0.5270378421266741
$ws.Cells.Item(9, 12).Value = 0.1700234219986214
$ws.Cells.Item(9, 13).Value = 0.2524592284763187
$ws.Cells.Item(10, 2).Value = 1.185427432860308
$ws.Cells.Item(10, 4).Value = 0.324740829660584
$ws.Cells.Item(10, 5).Value = 0.1749523125049448
$ws.Cells.Item(10, 6).Value = 5.189689608410447
$ws.Cells.Item(10, 7).Value = 0.002561216490893341
$ws.Cells.Item(10, 11).Value = 0.5871566670838888
$ws.Cells.Item(10, 12).Value = 0.1773129685504387
$ws.Cells.Item(10, 13).Value = 0.2598522962479137
$ws.Cells.Item(11, 2).Value = 1.198892094722964
$ws.Cells.Item(11, 4).Value = 0.3342720955928371
$ws.Cells.Item(11, 5).Value = 0.176690645013835
$ws.Cells.Item(11, 6).Value = 5.36750407090949
$ws.Cells.Item(11, 7).Value = 0.002557475449130158
$ws.Cells.Item(11, 11).Value = 0.6153256589580565
$ws.Cells.Item(11, 12).Value = 0.1807823210951085
$ws.Cells.Item(11, 13).Value = 0.2633910481458912
$ws.Cells.Item(12, 2).Value = 1.204090615221986
$ws.Cells.Item(12, 4).Value = 0.3378754782392548
$ws.Cells.Item(12, 5).Value = 0.1773419856612328
$ws.Cells.Item(12, 6).Value = 5.434938442781629
$ws.Cells.Item(12, 7).Value = 0.002556084199337149
$ws.Cells.Item(12, 11).Value = 0.6261126514941111
$ws.Cells.Item(12, 12).Value = 0.1821182523572134
$ws.Cells.Item(12, 13).Value = 0.2647563929908969
$ws.Cells.Item(13, 2).Value = 1.202966584273241
$ws.Cells.Item(13, 4).Value = 0.3370996806785058
$ws.Cells.Item(13, 5).Value = 0.177202012898273
$ws.Cells.Item(13, 6).Value = 5.420410711849001
$ws.Cells.Item(13, 7).Value = 0.002556382702562155
$ws.Cells.Item(13, 11).Value = 0.6237841100839034
$ws.Cells.Item(13, 12).Value = 0.181829548044135
$ws.Cells.Item(13, 13).Value = 0.2644612157576134
$ws.Cells.Item(14, 2).Value = 1.199317780967391
$ws.Cells.Item(14, 4).Value = 0.3345686637669871
$ws.Cells.Item(14, 5).Value = 0.1767443689605095
$ws.Cells.Item(14, 6).Value = 5.373049903045455
$ws.Cells.Item(14, 7).Value = 0.002557360482043372
$ws.Cells.Item(14, 11).Value = 0.6162106946287338
$ws.Cells.Item(14, 12).Value = 0.1808917841001971
$ws.Cells.Item(14, 13).Value = 0.2635028686500078
$ws.Cells.Item(15, 2).Value = 1.197095773983364
$ws.Cells.Item(15, 4).Value = 0.3330175875462658
$ws.Cells.Item(15, 5).Value = 0.176463152233457
$ws.Cells.Item(15, 6).Value = 5.344053215999736
$ws.Cells.Item(15, 7).Value = 0.002557962703700244
$ws.Cells.Item(15, 11).Value = 0.6115874473708516
$ws.Cells.Item(15, 12).Value = 0.1803202664181356
$ws.Cells.Item(15, 13).Value = 0.262919149055918
$ws.Cells.Item(16, 2).Value = 1.184561442392919
$ws.Cells.Item(16, 4).Value = 0.324117089106835
$ws.Cells.Item(16, 5).Value = 0.1748377285490266
$ws.Cells.Item(16, 6).Value = 5.178082688682792
$ws.Cells.Item(16, 7).Value = 0.002561464539443035
$ws.Cells.Item(16, 11).Value = 0.5853324406828335
$ws.Cells.Item(16, 12).Value = 0.1770893330439094
$ws.Cells.Item(16, 13).Value = 0.2596245678770615
$ws.Cells.Item(17, 2).Value = 1.177049675553349
$ws.Cells.Item(17, 4).Value = 0.318645937210448
$ws.Cells.Item(17, 5).Value = 0.1738280409085702
$ws.Cells.Item(17, 6).Value = 5.076437087288753
$ws.Cells.Item(17, 7).Value = 0.002563658207679259
$ws.Cells.Item(17, 11).Value = 0.5694374410355465
$ws.Cells.Item(17, 12).Value = 0.1751466072406203
$ws.Cells.Item(17, 13).Value = 0.2576484569709194
$ws.Cells.Item(18, 2).Value = 1.172794383523296
$ws.Cells.Item(18, 4).Value = 0.3154948967695788
$ws.Cells.Item(18, 5).Value = 0.1732426251487764
$ws.Cells.Item(18, 6).Value = 5.018034629293254
$ws.Cells.Item(18, 7).Value = 0.002564936682570882
$ws.Cells.Item(18, 11).Value = 0.5603721781120328
$ws.Cells.Item(18, 12).Value = 0.174043624909288
$ws.Cells.Item(18, 13).Value = 0.2565283765224109
$ws.Cells.Item(19, 2).Value = 1.171364825926247
$ws.Cells.Item(19, 4).Value = 0.3144272775376322
$ws.Cells.Item(19, 5).Value = 0.1730436056002347
$ws.Cells.Item(19, 6).Value = 4.998270938322491
$ws.Cells.Item(19, 7).Value = 0.002565372431076759
$ws.Cells.Item(19, 11).Value = 0.5573160210003607
$ws.Cells.Item(19, 12).Value = 0.1736726475168808
$ws.Cells.Item(19, 13).Value = 0.2561519736352835
$ws.Cells.Item(20, 2).Value = 1.177842559896476
$ws.Cells.Item(20, 4).Value = 0.3192287799456324
$ws.Cells.Item(20, 5).Value = 0.1739360058582369
$ws.Cells.Item(20, 6).Value = 5.08725103998097
$ws.Cells.Item(20, 7).Value = 0.002563422957091017
$ws.Cells.Item(20, 11).Value = 0.5711214948355234
$ws.Cells.Item(20, 12).Value = 0.1753519202345046
$ws.Cells.Item(20, 13).Value = 0.257857106633665
$ws.Cells.Item(21, 2).Value = 1.200386815718247
$ws.Cells.Item(21, 4).Value = 0.3353122413936944
$ws.Cells.Item(21, 5).Value = 0.1768789765893954
$ws.Cells.Item(21, 6).Value = 5.386958164738644
$ws.Cells.Item(21, 7).Value = 0.002557072596466896
$ws.Cells.Item(21, 11).Value = 0.618431917264445
$ws.Cells.Item(21, 12).Value = 0.1811666256523523
$ws.Cells.Item(21, 13).Value = 0.2637836715266459
$ws.Cells.Item(22, 2).Value = 1.215702195915156
$ws.Cells.Item(22, 4).Value = 0.3457895261798001
$ws.Cells.Item(22, 5).Value = 0.178762084735558
$ws.Cells.Item(22, 6).Value = 5.583419000085541
$ws.Cells.Item(22, 7).Value = 0.002553070257068012
$ws.Cells.Item(22, 11).Value = 0.6500524580296485
$ws.Cells.Item(22, 12).Value = 0.185096110536648
$ws.Cells.Item(22, 13).Value = 0.2678044929755785
$ws.Cells.Item(23, 2).Value = 1.207474883012765
$ws.Cells.Item(23, 4).Value = 0.3402005831257782
$ws.Cells.Item(23, 5).Value = 0.1777606583130433
$ws.Cells.Item(23, 6).Value = 5.478508774065858
$ws.Cells.Item(23, 7).Value = 0.002555192890091344
$ws.Cells.Item(23, 11).Value = 0.6331112192804653
$ws.Cells.Item(23, 12).Value = 0.1829870051971199
$ws.Cells.Item(23, 13).Value = 0.2656449962029583
$ws.Cells.Item(24, 2).Value = 1.177483899649502
$ws.Cells.Item(24, 4).Value = 0.3189652942282919
$ws.Cells.Item(24, 5).Value = 0.1738872102836426
$ws.Cells.Item(24, 6).Value = 5.082361943519942
$ws.Cells.Item(24, 7).Value = 0.002563529259943917
$ws.Cells.Item(24, 11).Value = 0.5703599069340441
$ws.Cells.Item(24, 12).Value = 0.1752590548915549
$ws.Cells.Item(24, 13).Value = 0.2577627262504478
$ws.Cells.Item(25, 2).Value = 1.147629313386005
$ws.Cells.Item(25, 4).Value = 0.2958980000023246
$ws.Cells.Item(25, 5).Value = 0.1695297586093947
$ws.Cells.Item(25, 6).Value = 4.657377551885247
$ws.Cells.Item(25, 7).Value = 0.002573161842663153
$ws.Cells.Item(25, 11).Value = 0.5056128248428422
$ws.Cells.Item(25, 12).Value = 0.1674741063597622
$ws.Cells.Item(25, 13).Value = 0.2498920570404657
